$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row renames
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# Title-case connector words (de/del/la/el/los/las/y) in state/municipality names
$ws.Range('B5').Value = 'Rincón De Romos'
$ws.Range('B24').Value = 'Amatenango De La Frontera'
$ws.Range('B27').Value = 'Benemérito De Las Américas'
$ws.Range('B34').Value = 'Comitán De Domínguez'
$ws.Range('B48').Value = 'Montecristo De Guerrero'
$ws.Range('B58').Value = 'Salto De Agua'
$ws.Range('B59').Value = 'San Cristóbal De Las Casas'
$ws.Range('B83').Value = 'Hidalgo Del Parral'
$ws.Range('A91').Value = 'Ciudad De México'
$ws.Range('A106').Value = 'Coahuila De Zaragoza'
$ws.Range('B116').Value = 'San Juan De Sabinas'
$ws.Range('B135').Value = 'San Pedro Del Gallo'
$ws.Range('A140').Value = 'Estado De México'
$ws.Range('B140').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B142').Value = 'Almoloya De Alquisiras'
$ws.Range('B144').Value = 'Atizapán De Zaragoza'
$ws.Range('B148').Value = 'Coacalco De Berriozábal'
$ws.Range('B152').Value = 'Ecatepec De Morelos'
$ws.Range('B158').Value = 'Naucalpan De Juárez'
$ws.Range('B162').Value = 'San Felipe Del Progreso'
$ws.Range('B168').Value = 'Tlalnepantla De Baz'
$ws.Range('B171').Value = 'Valle De Chalco Solidaridad'
$ws.Range('B175').Value = 'Apaseo El Alto'
$ws.Range('B176').Value = 'Apaseo El Grande'
$ws.Range('B185').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B195').Value = 'San Francisco Del Rincón'
$ws.Range('B197').Value = 'San Luis De La Paz'
$ws.Range('B198').Value = 'San Miguel De Allende'
$ws.Range('B199').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B201').Value = 'Valle De Santiago'
$ws.Range('B207').Value = 'Acapulco De Juárez'
$ws.Range('B210').Value = 'Atenango Del Río'
$ws.Range('B212').Value = 'Atoyac De Álvarez'
$ws.Range('B213').Value = 'Ayutla De Los Libres'
$ws.Range('B215').Value = 'Buenavista De Cuéllar'
$ws.Range('B216').Value = 'Chilapa De Álvarez'
$ws.Range('B217').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B219').Value = 'Coyuca De Benítez'
$ws.Range('B220').Value = 'Coyuca De Catalán'
$ws.Range('B222').Value = 'Cuetzala Del Progreso'
$ws.Range('B223').Value = 'Cutzamala De Pinzón'
$ws.Range('B227').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B228').Value = 'Iguala De La Independencia'
$ws.Range('B240').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B243').Value = 'Tlapa De Comonfort'
$ws.Range('B245').Value = 'Técpan De Galeana'
$ws.Range('B250').Value = 'Atotonilco El Grande'
$ws.Range('B259').Value = 'Huejutla De Reyes'
$ws.Range('B261').Value = 'Jacala De Ledezma'
$ws.Range('B266').Value = 'Nopala De Villagrán'
$ws.Range('B267').Value = 'Pachuca De Soto'
$ws.Range('B273').Value = 'Tenango De Doria'
$ws.Range('B275').Value = 'Tepehuacán De Guerrero'
$ws.Range('B276').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B277').Value = 'Tezontepec De Aldama'
$ws.Range('B283').Value = 'Tulancingo De Bravo'
$ws.Range('B285').Value = 'Zacualtipán De Ángeles'
$ws.Range('B286').Value = 'Zapotlán De Juárez'
$ws.Range('B290').Value = 'Autlán De Navarro'
$ws.Range('B294').Value = 'Encarnación De Díaz'
$ws.Range('B298').Value = 'Lagos De Moreno'
$ws.Range('B301').Value = 'Ojuelos De Jalisco'
$ws.Range('B304').Value = 'San Diego De Alejandría'
$ws.Range('B306').Value = 'Santa María De Los Ángeles'
$ws.Range('B307').Value = 'Santa María Del Oro'
$ws.Range('B309').Value = 'Talpa De Allende'
$ws.Range('B310').Value = 'Tamazula De Gordiano'
$ws.Range('B312').Value = 'Tepatitlán De Morelos'
$ws.Range('B314').Value = 'Tizapán El Alto'
$ws.Range('B316').Value = 'Yahualica De González Gallo'
$ws.Range('B319').Value = 'Zapotlán El Grande'
$ws.Range('A321').Value = 'Michoacán De Ocampo'
$ws.Range('B325').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B365').Value = 'Puente De Ixtla'
$ws.Range('B374').Value = 'Santa María Del Oro'
$ws.Range('B383').Value = 'Mier Y Noriega'
$ws.Range('B387').Value = 'San Nicolás De Los Garza'
$ws.Range('B389').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B394').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B395').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B396').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B397').Value = 'Huajuapan De León'
$ws.Range('B402').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B403').Value = 'Mártires De Tacubaya'
$ws.Range('B404').Value = 'Oaxaca De Juárez'
$ws.Range('B405').Value = 'Putla Villa De Guerrero'
$ws.Range('B423').Value = 'San Miguel El Grande'
$ws.Range('B430').Value = 'Santa Inés Del Monte'
$ws.Range('B433').Value = 'Santa María Jalapa Del Marqués'
$ws.Range('B444').Value = 'Tataltepec De Valdés'
$ws.Range('B445').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B460').Value = 'Izúcar De Matamoros'
$ws.Range('B465').Value = 'Los Reyes De Juárez'
$ws.Range('B472').Value = 'San Salvador El Seco'
$ws.Range('B475').Value = 'Tepanco De López'
$ws.Range('B477').Value = 'Tetela De Ocampo'
$ws.Range('B480').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B494').Value = 'Amealco De Bonfil'
$ws.Range('B499').Value = 'Jalpan De Serra'
$ws.Range('B500').Value = 'Landa De Matamoros'
$ws.Range('B502').Value = 'Pinal De Amoles'
$ws.Range('B504').Value = 'San Juan Del Río'
$ws.Range('B520').Value = 'Ciudad Del Maíz'
$ws.Range('B528').Value = 'Mexquitic De Carmona'
$ws.Range('B533').Value = 'San Ciro De Acosta'
$ws.Range('B539').Value = 'Santa María Del Río'
$ws.Range('B545').Value = 'Tanquián De Escobedo'
$ws.Range('B550').Value = 'Villa De Arriaga'
$ws.Range('B551').Value = 'Villa De Guadalupe'
$ws.Range('B552').Value = 'Villa De Ramos'
$ws.Range('B553').Value = 'Villa De Reyes'
$ws.Range('B600').Value = 'Soto La Marina'
$ws.Range('B610').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B612').Value = 'Tetla De La Solidaridad'
$ws.Range('A616').Value = 'Veracruz De Ignacio De La Llave'
$ws.Range('B619').Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range('B621').Value = 'Amatlán De Los Reyes'
$ws.Range('B627').Value = 'Boca Del Río'
$ws.Range('B628').Value = 'Camarón De Tejeda'
$ws.Range('B630').Value = 'Castillo De Teayo'
$ws.Range('B640').Value = 'Cosamaloapan De Carpio'
$ws.Range('B641').Value = 'Cosautlán De Carvajal'
$ws.Range('B655').Value = 'Ignacio De La Llave'
$ws.Range('B667').Value = 'Las Vigas De Ramírez'
$ws.Range('B668').Value = 'Lerdo De Tejada'
$ws.Range('B670').Value = 'Martínez De La Torre'
$ws.Range('B671').Value = 'Medellín De Bravo'
$ws.Range('B676').Value = 'Ozuluama De Mascareñas'
$ws.Range('B678').Value = 'Paso Del Macho'
$ws.Range('B680').Value = 'Poza Rica De Hidalgo'
$ws.Range('B685').Value = 'Sayula De Alemán'
$ws.Range('B687').Value = 'Soledad De Doblado'
$ws.Range('B691').Value = 'Tatahuicapan De Juárez'
$ws.Range('B732').Value = 'Nochistlán De Mejía'
$ws.Range('B739').Value = 'Teúl De González Ortega'
$ws.Range('B742').Value = 'Villa De Cos'

# TOTAL -> Total
$ws.Range('A746').Value = 'Total'

# Remove trailing footnote rows (749-753)
$ws.Range("A749:A753").EntireRow.Delete()
